$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 6057.4165
$ws.Range("J32").Value = 3402.25
$ws.Range("L32").Value = 3402.25
$ws.Range("N32").Value = -4054.25
$ws.Range("H53").Value = 1679.8
$ws.Range("I53").Value = 333.6
$ws.Range("J53").Value = 2352.9
$ws.Range("K53").Value = 333.6
$ws.Range("L53").Value = 2352.9
$ws.Range("M53").Value = 303.4
$ws.Range("N53").Value = -3626.9
$ws.Range("H86").Value = 2814.8462
$ws.Range("I86").Value = 2580.95
$ws.Range("J86").Value = 3061.0527
$ws.Range("K86").Value = 2580.95
$ws.Range("L86").Value = 3061.0527
$ws.Range("M86").Value = -1457.95
$ws.Range("N86").Value = -5307.0527
$ws.Range("H89").Value = 2814.8462
$ws.Range("I89").Value = 2580.95
$ws.Range("J89").Value = 3061.0527
$ws.Range("K89").Value = 12904.75
$ws.Range("L89").Value = 15305.2635
$ws.Range("M89").Value = -7288.75
$ws.Range("N89").Value = -26537.2635
$ws.Range("H99").Value = 338.33334
$ws.Range("I99").Value = 338.33334
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1015.00002
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 482.9999799999999
$ws.Range("N99").ClearContents()
$ws.Range("H106").Value = 5368.7144
$ws.Range("I106").Value = 4017.2222
$ws.Range("J106").Value = 7801.4
$ws.Range("K106").Value = 4017.2222
$ws.Range("L106").Value = 7801.4
$ws.Range("M106").Value = -3386.2222
$ws.Range("N106").Value = -9063.4
$ws.Range("H112").Value = 1290.48
$ws.Range("J112").Value = 1419.2
$ws.Range("L112").Value = 4257.6
$ws.Range("N112").Value = -6473.6
$ws.Range("H113").Value = 7209.15
$ws.Range("I113").Value = 13797.6
$ws.Range("K113").Value = 13797.6
$ws.Range("M113").Value = -10543.6
$ws.Range("H115").Value = 352
$ws.Range("I115").Value = 352
$ws.Range("K115").Value = 1056
$ws.Range("M115").Value = 511
$ws.Range("H132").Value = 2761
$ws.Range("I132").Value = 2717.611
$ws.Range("K132").Value = 8152.833
$ws.Range("M132").Value = -5622.833

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3111.8655
$ws.Range("I32").Value = 2475.5715
$ws.Range("J32").Value = 13504.667
$ws.Range("K32").Value = 2475.5715
$ws.Range("L32").Value = 13504.667
$ws.Range("M32").Value = -2188.5715
$ws.Range("N32").Value = -14078.667
$ws.Range("H45").Value = 2050.1052
$ws.Range("I45").Value = 2067.6875
$ws.Range("K45").Value = 2067.6875
$ws.Range("M45").Value = -1690.6875
$ws.Range("H74").Value = 17547954
$ws.Range("I74").Value = 27780132
$ws.Range("J74").Value = 7076.5713
$ws.Range("K74").Value = 27780132
$ws.Range("L74").Value = 7076.5713
$ws.Range("M74").Value = -27779258
$ws.Range("N74").Value = -8824.5713
$ws.Range("H76").Value = 599288
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 599288
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 599288
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -599964
$ws.Range("H77").Value = 17547954
$ws.Range("I77").Value = 27780132
$ws.Range("J77").Value = 7076.5713
$ws.Range("K77").Value = 138900660
$ws.Range("L77").Value = 35382.85649999999
$ws.Range("M77").Value = -138896292
$ws.Range("N77").Value = -44118.85649999999
$ws.Range("H79").Value = 599288
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 599288
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 599288
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -601628
$ws.Range("H102").Value = 3702.6
$ws.Range("I102").Value = 3273.9
$ws.Range("J102").Value = 4560
$ws.Range("K102").Value = 3273.9
$ws.Range("L102").Value = 4560
$ws.Range("M102").Value = -1651.9
$ws.Range("N102").Value = -7804

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3914.5293
$ws.Range("I20").Value = 3571.2727
$ws.Range("K20").Value = 3571.2727
$ws.Range("M20").Value = -3324.2727
$ws.Range("H21").Value = 56653.332
$ws.Range("J21").Value = 56653.332
$ws.Range("L21").Value = 56653.332
$ws.Range("N21").Value = -57125.332
$ws.Range("H107").Value = 2698.6667
$ws.Range("I107").Value = 2349.5
$ws.Range("J107").Value = 2873.25
$ws.Range("K107").Value = 2349.5
$ws.Range("L107").Value = 2873.25
$ws.Range("M107").Value = -429.5
$ws.Range("N107").Value = -6713.25
$ws.Range("H140").Value = 51612.777
$ws.Range("J140").Value = 51612.777
$ws.Range("L140").Value = 51612.777
$ws.Range("N140").Value = -61972.777

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4944.7915
$ws.Range("I58").Value = 1981.6111
$ws.Range("K58").Value = 1981.6111
$ws.Range("M58").Value = -1778.6111
$ws.Range("H111").Value = 73750.39999999999
$ws.Range("J111").Value = 73750.39999999999
$ws.Range("L111").Value = 73750.39999999999
$ws.Range("N111").Value = -81930.39999999999
$ws.Range("H134").Value = 3341.5173
$ws.Range("I134").Value = 2473.9565
$ws.Range("J134").Value = 6667.1665
$ws.Range("K134").Value = 7421.869499999999
$ws.Range("L134").Value = 20001.4995
$ws.Range("M134").Value = -4886.869499999999
$ws.Range("N134").Value = -25071.4995
$ws.Range("H136").Value = 4944.7915
$ws.Range("I136").Value = 1981.6111
$ws.Range("K136").Value = 5944.8333
$ws.Range("M136").Value = -3394.8333
$ws.Range("H141").Value = 331054
$ws.Range("J141").Value = 331054
$ws.Range("L141").Value = 331054
$ws.Range("N141").Value = -341414

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1318.65
$ws.Range("I98").Value = 1004.8889
$ws.Range("J98").Value = 1575.3636
$ws.Range("K98").Value = 3014.6667
$ws.Range("L98").Value = 4726.0908
$ws.Range("M98").Value = -1516.6667
$ws.Range("N98").Value = -7722.0908
$ws.Range("H112").Value = 50003250
$ws.Range("I112").Value = 166668340
$ws.Range("K112").Value = 500005020
$ws.Range("M112").Value = -500003912

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1916.0714
$ws.Range("J97").Value = 2286.5715
$ws.Range("L97").Value = 2286.5715
$ws.Range("N97").Value = -3278.5715
$ws.Range("H122").Value = 7627.385
$ws.Range("I122").Value = 6915.6
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 20746.8
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -18296.8
$ws.Range("N122").Value = -34900

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3308.1428
$ws.Range("I68").Value = 3443.4443
$ws.Range("J68").Value = 2496.3333
$ws.Range("K68").Value = 3443.4443
$ws.Range("L68").Value = 2496.3333
$ws.Range("M68").Value = -2694.4443
$ws.Range("N68").Value = -3994.3333
$ws.Range("H71").Value = 3308.1428
$ws.Range("I71").Value = 3443.4443
$ws.Range("J71").Value = 2496.3333
$ws.Range("K71").Value = 17217.2215
$ws.Range("L71").Value = 12481.6665
$ws.Range("M71").Value = -13473.2215
$ws.Range("N71").Value = -19969.6665
$ws.Range("H93").Value = 2000
$ws.Range("J93").Value = 2000
$ws.Range("L93").Value = 2000
$ws.Range("N93").Value = -4496
$ws.Range("H97").Value = 14500
$ws.Range("I97").Value = 9000
$ws.Range("J97").Value = 20000
$ws.Range("K97").Value = 9000
$ws.Range("L97").Value = 20000
$ws.Range("M97").Value = -8009
$ws.Range("N97").Value = -21982
$ws.Range("H132").Value = 9130.566000000001
$ws.Range("I132").Value = 7806.8335
$ws.Range("K132").Value = 23420.5005
$ws.Range("M132").Value = -20890.5005

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 11901
$ws.Range("I30").Value = 11901
$ws.Range("K30").Value = 11901
$ws.Range("M30").Value = -11794
$ws.Range("H62").Value = 8484.308000000001
$ws.Range("J62").Value = 11949.5
$ws.Range("L62").Value = 11949.5
$ws.Range("N62").Value = -13197.5
$ws.Range("H65").Value = 8484.308000000001
$ws.Range("J65").Value = 11949.5
$ws.Range("L65").Value = 59747.5
$ws.Range("N65").Value = -65987.5
$ws.Range("H94").Value = 14000
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H113").Value = 390.21817
$ws.Range("I113").Value = 236.17949
$ws.Range("J113").Value = 765.6875
$ws.Range("K113").Value = 708.53847
$ws.Range("L113").Value = 2297.0625
$ws.Range("M113").Value = 1461.46153
$ws.Range("N113").Value = -6637.0625
$ws.Range("H132").Value = 5860.722
$ws.Range("I132").Value = 3713.4285
$ws.Range("K132").Value = 11140.2855
$ws.Range("M132").Value = -8610.2855
